# ===================================================================
# Sections sheet: rebuild the Sections table, add new course sections
# ===================================================================
$wb = $excel.ActiveWorkbook
$wsSections = $wb.Worksheets.Item("Sections")
$ws = $wsSections

# Header: "Type" -> "Section Number"
$ws.Range("B1").Value = "Section Number"

# Existing rows 2 and 3 currently hold trailing-comma text values;
# clear them first so the re-written cells pick up the column's
# inherited formatting (matches the newly added rows below).
$ws.Range("A2:F3").ClearContents()

# Row 2: MATHF111 L1
$ws.Range("A2").Value = "MATHF111"
$ws.Range("B2").Value = "L1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "M,W"
$ws.Range("E2").Value = 5102
$ws.Range("F2").Value = "Trilok,Keskar"
$ws.Rows.Item(2).RowHeight = 15

# Row 3: MATHF111 L2
$ws.Range("A3").Value = "MATHF111"
$ws.Range("B3").Value = "L2"
$ws.Range("C3").Value = 34
$ws.Range("D3").Value = "T,Th,S"
$ws.Range("E3").Value = 5105
$ws.Range("F3").Value = "Anirudha,Sharma"
$ws.Rows.Item(3).RowHeight = 15

# Row 4: CSF111 T1
$ws.Range("A4").Value = "CSF111"
$ws.Range("B4").Value = "T1"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "W"
$ws.Range("E4").Value = 6102
$ws.Range("F4").Value = "Vinti"
$ws.Rows.Item(4).RowHeight = 15

# Row 5: CSF111 T2
$ws.Range("A5").Value = "CSF111"
$ws.Range("B5").Value = "T2"
$ws.Range("C5").Value = 89
$ws.Range("D5").Value = "S"
$ws.Range("E5").Value = 2104
$ws.Range("F5").Value = "Prakhar"
$ws.Rows.Item(5).RowHeight = 15

# Row 6: CHEMF111 L1
$ws.Range("A6").Value = "CHEMF111"
$ws.Range("B6").Value = "L1"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "M,T,F"
$ws.Range("E6").Value = 5105
$ws.Range("F6").Value = "Daksh Jain,Ram Srivastava"
$ws.Rows.Item(6).RowHeight = 15

# Row 7: CHEMF111 L2
$ws.Range("A7").Value = "CHEMF111"
$ws.Range("B7").Value = "L2"
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = "T,Th,S"
$ws.Range("E7").Value = 5102
$ws.Range("F7").Value = "Dhyann,Prakhar"
$ws.Rows.Item(7).RowHeight = 15

# Row 8: CHEMF111 T3
$ws.Range("A8").Value = "CHEMF111"
$ws.Range("B8").Value = "T3"
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = "F"
$ws.Range("E8").Value = 6164
$ws.Range("F8").Value = "Paritosh"
$ws.Rows.Item(8).RowHeight = 15

# Row 9: CSF111 T3
$ws.Range("A9").Value = "CSF111"
$ws.Range("B9").Value = "T3"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "W"
$ws.Range("E9").Value = 6168
$ws.Range("F9").Value = "Ram"
$ws.Rows.Item(9).RowHeight = 15

# Row 10: CSF111 L1
$ws.Range("A10").Value = "CSF111"
$ws.Range("B10").Value = "L1"
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = "Th"
$ws.Range("E10").Value = 6109
$ws.Range("F10").Value = "Abhishek"
$ws.Rows.Item(10).RowHeight = 15

# Row 11: MATHF111 T1
$ws.Range("A11").Value = "MATHF111"
$ws.Range("B11").Value = "T1"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Th"
$ws.Range("E11").Value = 6168
$ws.Range("F11").Value = "Divyum"

# Column widths: widen the professor-list column and give the new
# Code/Section-Number columns an explicit width
$ws.Range("A1:B11").ColumnWidth = 13.5
$ws.Columns.Item(6).ColumnWidth = 22.33

# ===================================================================
# Courses sheet: de-duplicate the redundant date/time number format
# (C4 used a second, equivalent numFmt -- align it with C3's format)
# ===================================================================
$wsCourses = $wb.Worksheets.Item("Courses")
$wsCourses.Range("C3").Copy()
$wsCourses.Range("C4").PasteSpecial(-4122)

# Sections stays the active sheet/selection, matching the original
# workbook's active-tab state
$ws.Activate()
$ws.Range("E11").Select()
